## Dispatch Priority by Elec Source - add three new fuel rows (crude oil,
## heavy/residual fuel oil, municipal solid waste) whose priority values
## mirror existing rows, plus a bold/wrapped title for column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# --- New data rows (15-17). Column B gets a single formula copying the
# priority value from an existing row; C:AK are then filled with the same
# relative formula (mirrors the pattern already used by rows 13 & 14). ---

# Row 15: crude oil -> mirrors "petroleum" (row 11)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

# Row 16: heavy or residual fuel oil -> mirrors "petroleum" (row 11)
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

# Row 17: municipal solid waste -> mirrors "biomass" (row 9)
$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# --- Column A header / title cell, bold + wrapped, with a taller row 1 ---
$ws.Range("A1").Value = "Dispatch Priority (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Widen column A slightly to fit the new wrapped header text
$ws.Columns.Item(1).ColumnWidth = 22.93
